# --- NATMI TPM refresh -------------------------------------------------------
# The commit recomputes the NATMI ligand/receptor (Tnfsf11 -> Tnfrsf11a) scores
# with an updated TPM table. Net effect on this worksheet:
#   1) The "Inflammatory-Mac" sending-cluster rows are dropped entirely, and the
#      "MuSCs" sending-cluster rows that followed them shift up to take their place
#      (the shared-string table reshuffle visible in the raw XML diff is just the
#      automatic side effect of this row deletion - Excel/COM manages that table).
#   2) Every remaining row gets refreshed values in columns E:T.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Inflammatory-Mac" sending-cluster block (rows 10-13); the cells
# below shift up so the "MuSCs" block (previously rows 14-17) becomes rows 10-13.
$ws.Range("A10:T13").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# Row 2: Sending cluster = ECs; Target cluster = ECs
$ws.Cells.Item(2, 5).Value = 1  # E2
$ws.Cells.Item(2, 6).Value = 0.3333333333333333  # F2
$ws.Cells.Item(2, 7).Value = 0.05344666666666667  # G2
$ws.Cells.Item(2, 8).Value = 0.16034  # H2
$ws.Cells.Item(2, 9).Value = 0.01683020325561672  # I2
$ws.Cells.Item(2, 10).Value = 0.01683020325561672  # J2
$ws.Cells.Item(2, 11).Value = 3  # K2
$ws.Cells.Item(2, 12).Value = 1  # L2
$ws.Cells.Item(2, 13).Value = 2.386127333333333  # M2
$ws.Cells.Item(2, 14).Value = 7.158382  # N2
$ws.Cells.Item(2, 15).Value = 0.04731171100836582  # O2
$ws.Cells.Item(2, 16).Value = 0.04731171100836583  # P2
$ws.Cells.Item(2, 17).Value = 0.1275305522088889  # Q2
$ws.Cells.Item(2, 18).Value = 1.14777496988  # R2
$ws.Cells.Item(2, 19).Value = 0.0007962657126417957  # S2
$ws.Cells.Item(2, 20).Value = 0.0007962657126417958  # T2

# Row 3: Sending cluster = ECs; Target cluster = Inflammatory-Mac
$ws.Cells.Item(3, 5).Value = 1  # E3
$ws.Cells.Item(3, 6).Value = 0.3333333333333333  # F3
$ws.Cells.Item(3, 7).Value = 0.05344666666666667  # G3
$ws.Cells.Item(3, 8).Value = 0.16034  # H3
$ws.Cells.Item(3, 9).Value = 0.01683020325561672  # I3
$ws.Cells.Item(3, 10).Value = 0.01683020325561672  # J3
$ws.Cells.Item(3, 11).Value = 3  # K3
$ws.Cells.Item(3, 12).Value = 1  # L3
$ws.Cells.Item(3, 13).Value = 26.214127  # M3
$ws.Cells.Item(3, 14).Value = 78.642381  # N3
$ws.Cells.Item(3, 15).Value = 0.519769076710603  # O3
$ws.Cells.Item(3, 16).Value = 0.519769076710603  # P3
$ws.Cells.Item(3, 17).Value = 1.401057707726667  # Q3
$ws.Cells.Item(3, 18).Value = 12.60951936954  # R3
$ws.Cells.Item(3, 19).Value = 0.008747819207023687  # S3
$ws.Cells.Item(3, 20).Value = 0.008747819207023687  # T3

# Row 4: Sending cluster = ECs; Target cluster = MuSCs
$ws.Cells.Item(4, 5).Value = 1  # E4
$ws.Cells.Item(4, 6).Value = 0.3333333333333333  # F4
$ws.Cells.Item(4, 7).Value = 0.05344666666666667  # G4
$ws.Cells.Item(4, 8).Value = 0.16034  # H4
$ws.Cells.Item(4, 9).Value = 0.01683020325561672  # I4
$ws.Cells.Item(4, 10).Value = 0.01683020325561672  # J4
$ws.Cells.Item(4, 11).Value = 1  # K4
$ws.Cells.Item(4, 12).Value = 0.3333333333333333  # L4
$ws.Cells.Item(4, 13).Value = 0.007255999999999999  # M4
$ws.Cells.Item(4, 14).Value = 0.021768  # N4
$ws.Cells.Item(4, 15).Value = 0.0001438706854747494  # O4
$ws.Cells.Item(4, 16).Value = 0.0001438706854747494  # P4
$ws.Cells.Item(4, 17).Value = 0.0003878090133333333  # Q4
$ws.Cells.Item(4, 18).Value = 0.00349028112  # R4
$ws.Cells.Item(4, 19).Value = 0.000002421372879064935  # S4
$ws.Cells.Item(4, 20).Value = 0.000002421372879064936  # T4

# Row 5: Sending cluster = ECs; Target cluster = Resolving-Mac
$ws.Cells.Item(5, 5).Value = 1  # E5
$ws.Cells.Item(5, 6).Value = 0.3333333333333333  # F5
$ws.Cells.Item(5, 7).Value = 0.05344666666666667  # G5
$ws.Cells.Item(5, 8).Value = 0.16034  # H5
$ws.Cells.Item(5, 9).Value = 0.01683020325561672  # I5
$ws.Cells.Item(5, 10).Value = 0.01683020325561672  # J5
$ws.Cells.Item(5, 11).Value = 3  # K5
$ws.Cells.Item(5, 12).Value = 1  # L5
$ws.Cells.Item(5, 13).Value = 21.82666933333333  # M5
$ws.Cells.Item(5, 14).Value = 65.480008  # N5
$ws.Cells.Item(5, 15).Value = 0.4327753415955564  # O5
$ws.Cells.Item(5, 16).Value = 0.4327753415955564  # P5
$ws.Cells.Item(5, 17).Value = 1.166562720302222  # Q5
$ws.Cells.Item(5, 18).Value = 10.49906448272  # R5
$ws.Cells.Item(5, 19).Value = 0.007283696963072171  # S5
$ws.Cells.Item(5, 20).Value = 0.007283696963072171  # T5

# Row 6: Sending cluster = FAPs; Target cluster = ECs
$ws.Cells.Item(6, 5).Value = 3  # E6
$ws.Cells.Item(6, 6).Value = 1  # F6
$ws.Cells.Item(6, 7).Value = 3.017399  # G6
$ws.Cells.Item(6, 8).Value = 9.052197  # H6
$ws.Cells.Item(6, 9).Value = 0.9501703593606328  # I6
$ws.Cells.Item(6, 10).Value = 0.9501703593606329  # J6
$ws.Cells.Item(6, 11).Value = 3  # K6
$ws.Cells.Item(6, 12).Value = 1  # L6
$ws.Cells.Item(6, 13).Value = 2.386127333333333  # M6
$ws.Cells.Item(6, 14).Value = 7.158382  # N6
$ws.Cells.Item(6, 15).Value = 0.04731171100836582  # O6
$ws.Cells.Item(6, 16).Value = 0.04731171100836583  # P6
$ws.Cells.Item(6, 17).Value = 7.199898229472665  # Q6
$ws.Cells.Item(6, 18).Value = 64.799084065254  # R6
$ws.Cells.Item(6, 19).Value = 0.04495418545078535  # S6
$ws.Cells.Item(6, 20).Value = 0.04495418545078537  # T6

# Row 7: Sending cluster = FAPs; Target cluster = Inflammatory-Mac
$ws.Cells.Item(7, 5).Value = 3  # E7
$ws.Cells.Item(7, 6).Value = 1  # F7
$ws.Cells.Item(7, 7).Value = 3.017399  # G7
$ws.Cells.Item(7, 8).Value = 9.052197  # H7
$ws.Cells.Item(7, 9).Value = 0.9501703593606328  # I7
$ws.Cells.Item(7, 10).Value = 0.9501703593606329  # J7
$ws.Cells.Item(7, 11).Value = 3  # K7
$ws.Cells.Item(7, 12).Value = 1  # L7
$ws.Cells.Item(7, 13).Value = 26.214127  # M7
$ws.Cells.Item(7, 14).Value = 78.642381  # N7
$ws.Cells.Item(7, 15).Value = 0.519769076710603  # O7
$ws.Cells.Item(7, 16).Value = 0.519769076710603  # P7
$ws.Cells.Item(7, 17).Value = 79.098480595673  # Q7
$ws.Cells.Item(7, 18).Value = 711.886325361057  # R7
$ws.Cells.Item(7, 19).Value = 0.493869170402658  # S7
$ws.Cells.Item(7, 20).Value = 0.493869170402658  # T7

# Row 8: Sending cluster = FAPs; Target cluster = MuSCs
$ws.Cells.Item(8, 5).Value = 3  # E8
$ws.Cells.Item(8, 6).Value = 1  # F8
$ws.Cells.Item(8, 7).Value = 3.017399  # G8
$ws.Cells.Item(8, 8).Value = 9.052197  # H8
$ws.Cells.Item(8, 9).Value = 0.9501703593606328  # I8
$ws.Cells.Item(8, 10).Value = 0.9501703593606329  # J8
$ws.Cells.Item(8, 11).Value = 1  # K8
$ws.Cells.Item(8, 12).Value = 0.3333333333333333  # L8
$ws.Cells.Item(8, 13).Value = 0.007255999999999999  # M8
$ws.Cells.Item(8, 14).Value = 0.021768  # N8
$ws.Cells.Item(8, 15).Value = 0.0001438706854747494  # O8
$ws.Cells.Item(8, 16).Value = 0.0001438706854747494  # P8
$ws.Cells.Item(8, 17).Value = 0.02189424714399999  # Q8
$ws.Cells.Item(8, 18).Value = 0.197048224296  # R8
$ws.Cells.Item(8, 19).Value = 0.0001367016609190032  # S8
$ws.Cells.Item(8, 20).Value = 0.0001367016609190032  # T8

# Row 9: Sending cluster = FAPs; Target cluster = Resolving-Mac
$ws.Cells.Item(9, 5).Value = 3  # E9
$ws.Cells.Item(9, 6).Value = 1  # F9
$ws.Cells.Item(9, 7).Value = 3.017399  # G9
$ws.Cells.Item(9, 8).Value = 9.052197  # H9
$ws.Cells.Item(9, 9).Value = 0.9501703593606328  # I9
$ws.Cells.Item(9, 10).Value = 0.9501703593606329  # J9
$ws.Cells.Item(9, 11).Value = 3  # K9
$ws.Cells.Item(9, 12).Value = 1  # L9
$ws.Cells.Item(9, 13).Value = 21.82666933333333  # M9
$ws.Cells.Item(9, 14).Value = 65.480008  # N9
$ws.Cells.Item(9, 15).Value = 0.4327753415955564  # O9
$ws.Cells.Item(9, 16).Value = 0.4327753415955564  # P9
$ws.Cells.Item(9, 17).Value = 65.85977021973065  # Q9
$ws.Cells.Item(9, 18).Value = 592.737931977576  # R9
$ws.Cells.Item(9, 19).Value = 0.4112103018462704  # S9
$ws.Cells.Item(9, 20).Value = 0.4112103018462705  # T9

# Row 10: Sending cluster = MuSCs; Target cluster = ECs
$ws.Cells.Item(10, 5).Value = 1  # E10
$ws.Cells.Item(10, 6).Value = 0.3333333333333333  # F10
$ws.Cells.Item(10, 7).Value = 0.1047943333333333  # G10
$ws.Cells.Item(10, 8).Value = 0.314383  # H10
$ws.Cells.Item(10, 9).Value = 0.03299943738375047  # I10
$ws.Cells.Item(10, 10).Value = 0.03299943738375048  # J10
$ws.Cells.Item(10, 11).Value = 3  # K10
$ws.Cells.Item(10, 12).Value = 1  # L10
$ws.Cells.Item(10, 13).Value = 2.386127333333333  # M10
$ws.Cells.Item(10, 14).Value = 7.158382  # N10
$ws.Cells.Item(10, 15).Value = 0.04731171100836582  # O10
$ws.Cells.Item(10, 16).Value = 0.04731171100836583  # P10
$ws.Cells.Item(10, 17).Value = 0.2500526231451111  # Q10
$ws.Cells.Item(10, 18).Value = 2.250473608306  # R10
$ws.Cells.Item(10, 19).Value = 0.001561259844938666  # S10
$ws.Cells.Item(10, 20).Value = 0.001561259844938666  # T10

# Row 11: Sending cluster = MuSCs; Target cluster = Inflammatory-Mac
$ws.Cells.Item(11, 5).Value = 1  # E11
$ws.Cells.Item(11, 6).Value = 0.3333333333333333  # F11
$ws.Cells.Item(11, 7).Value = 0.1047943333333333  # G11
$ws.Cells.Item(11, 8).Value = 0.314383  # H11
$ws.Cells.Item(11, 9).Value = 0.03299943738375047  # I11
$ws.Cells.Item(11, 10).Value = 0.03299943738375048  # J11
$ws.Cells.Item(11, 11).Value = 3  # K11
$ws.Cells.Item(11, 12).Value = 1  # L11
$ws.Cells.Item(11, 13).Value = 26.214127  # M11
$ws.Cells.Item(11, 14).Value = 78.642381  # N11
$ws.Cells.Item(11, 15).Value = 0.519769076710603  # O11
$ws.Cells.Item(11, 16).Value = 0.519769076710603  # P11
$ws.Cells.Item(11, 17).Value = 2.747091962880333  # Q11
$ws.Cells.Item(11, 18).Value = 24.723827665923  # R11
$ws.Cells.Item(11, 19).Value = 0.01715208710092134  # S11
$ws.Cells.Item(11, 20).Value = 0.01715208710092134  # T11

# Row 12: Sending cluster = MuSCs; Target cluster = MuSCs
$ws.Cells.Item(12, 5).Value = 1  # E12
$ws.Cells.Item(12, 6).Value = 0.3333333333333333  # F12
$ws.Cells.Item(12, 7).Value = 0.1047943333333333  # G12
$ws.Cells.Item(12, 8).Value = 0.314383  # H12
$ws.Cells.Item(12, 9).Value = 0.03299943738375047  # I12
$ws.Cells.Item(12, 10).Value = 0.03299943738375048  # J12
$ws.Cells.Item(12, 11).Value = 1  # K12
$ws.Cells.Item(12, 12).Value = 0.3333333333333333  # L12
$ws.Cells.Item(12, 13).Value = 0.007255999999999999  # M12
$ws.Cells.Item(12, 14).Value = 0.021768  # N12
$ws.Cells.Item(12, 15).Value = 0.0001438706854747494  # O12
$ws.Cells.Item(12, 16).Value = 0.0001438706854747494  # P12
$ws.Cells.Item(12, 17).Value = 0.0007603876826666666  # Q12
$ws.Cells.Item(12, 18).Value = 0.006843489144000001  # R12
$ws.Cells.Item(12, 19).Value = 0.000004747651676681249  # S12
$ws.Cells.Item(12, 20).Value = 0.000004747651676681252  # T12

# Row 13: Sending cluster = MuSCs; Target cluster = Resolving-Mac
$ws.Cells.Item(13, 5).Value = 1  # E13
$ws.Cells.Item(13, 6).Value = 0.3333333333333333  # F13
$ws.Cells.Item(13, 7).Value = 0.1047943333333333  # G13
$ws.Cells.Item(13, 8).Value = 0.314383  # H13
$ws.Cells.Item(13, 9).Value = 0.03299943738375047  # I13
$ws.Cells.Item(13, 10).Value = 0.03299943738375048  # J13
$ws.Cells.Item(13, 11).Value = 3  # K13
$ws.Cells.Item(13, 12).Value = 1  # L13
$ws.Cells.Item(13, 13).Value = 21.82666933333333  # M13
$ws.Cells.Item(13, 14).Value = 65.480008  # N13
$ws.Cells.Item(13, 15).Value = 0.4327753415955564  # O13
$ws.Cells.Item(13, 16).Value = 0.4327753415955564  # P13
$ws.Cells.Item(13, 17).Value = 2.287311261673778  # Q13
$ws.Cells.Item(13, 18).Value = 20.585801355064  # R13
$ws.Cells.Item(13, 19).Value = 0.01428134278621378  # S13
$ws.Cells.Item(13, 20).Value = 0.01428134278621379  # T13
